# Distribution field add for excel upload in all modules
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I (existing "Cogm per kg" column shifts to J)
$ws.Columns("I").Insert()

# New column I: "Distribution channel code"
$ws.Range("I1").Value = "Distribution channel code"
$ws.Range("I2").Value = "TR"
$ws.Range("I3").Value = "GO"

# Match the style of the header row / match the formatting used elsewhere
$ws.Range("I1").Font.Bold = $true

# Column width adjustment: new column I gets a (non bestFit) custom width
$ws.Columns("I").ColumnWidth = 21.7
